$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "{d.records[i].date}"
$ws.Range("A3").Value = "{d.records[i+1].date}"
$ws.Range("B2").Value = "{d.records[i].test}"
$ws.Range("B3").Value = "{d.records[i+1].test}"
$ws.Range("D2").Value = "{d.records[i].status}"
$ws.Range("D3").Value = "{d.records[i+1].status}"

$ws.Range("D2").Select()
